$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 so it becomes the active cell (matches sheetView selection in diff)
$ws.Range("E8").Select()
